$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.134.41'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +6.72%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.102.74'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.62%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.13'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.095.19'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.64%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.528'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.150'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +13.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.75'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +8.44%  '
$ws.Range("E12").Value = '  +3.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000247'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +8.56%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.47'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.77%  '
$ws.Range("E15").Value = '  +0.81%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.616.28'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.16'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.52%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.038.63'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +6.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.101.29'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '464.55'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.14'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.74%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.725'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.60%  '
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.81'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.76%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.39'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.22'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.32%  '
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("E30").Value = '  +5.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.82'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +9.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.88'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.94%  '
$ws.Range("E33").Value = '  +3.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0861'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +13.44%  '
$ws.Range("E35").Value = '  +16.64%  '
$ws.Range("E36").Value = '  +5.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.31'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +20.50%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.03'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.96%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '50.82'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.73%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '432.40'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +9.00%  '
$ws.Range("E41").Value = '  +0.75%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.904.24'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.81%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0368'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.280'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +12.55%  '
$ws.Range("E45").Value = '  +6.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.15'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '35.25'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.08%  '
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '122.87'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.19%  '
$ws.Range("E50").Value = '  +1.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.45'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.29%  '
